$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Extend header row (row 1) with two new columns P and Q, matching the
# existing header style (copy format from O1, which already carries the
# bold/border header style).
$ws.Range("O1").Copy()
$ws.Range("P1:Q1").PasteSpecial(-4122)
$ws.Range("P1").Value = 14
$ws.Range("Q1").Value = 15

# For each data row, swap the values in columns I/K and M/O, then add the
# new P and Q columns with value 2.
for ($r = 2; $r -le 25; $r++) {
    $i = $ws.Cells.Item($r, 9).Value()   # I
    $k = $ws.Cells.Item($r, 11).Value()  # K
    $m = $ws.Cells.Item($r, 13).Value()  # M
    $o = $ws.Cells.Item($r, 15).Value()  # O

    $ws.Cells.Item($r, 9).Value = $k    # I <- old K
    $ws.Cells.Item($r, 11).Value = $i   # K <- old I
    $ws.Cells.Item($r, 13).Value = $o   # M <- old O
    $ws.Cells.Item($r, 15).Value = $m   # O <- old M

    $ws.Cells.Item($r, 16).Value = 2    # P
    $ws.Cells.Item($r, 17).Value = 2    # Q
}
